$d = $word.ActiveDocument

# Header/footer distance collapsed to 0 in the restored section properties.
$ps = $d.PageSetup
$ps.HeaderDistance = 0
$ps.FooterDistance = 0

# Merge the separate "Date" and ":" runs into a single "Date:" run.
$found = $d.Content.Find.Execute("Date:", $false, $false, $false, $false, $false, `
                                  $true, 1, $false, "Date:", 2)

# Append the new "Person" / "Message" timing block at the end of the document.
$lines = @(
    "Person:",
    "Total run time fan10= 0.8604384049976943 seconds.",
    "Total run time fan200= 0.6758918440027628 seconds.",
    "Message:",
    "Total run time fan10= 22.25640131099499 seconds.",
    "Total run time fan200= 14.956734024002799 seconds."
)

foreach ($line in $lines) {
    $p = $d.Paragraphs.Last
    $p.Range.InsertParagraphAfter()
    $newp = $d.Paragraphs.Last
    $newp.Range.Text = $line
}
